$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert 25 new worker rows above the existing last data row (row 17),
#    which pushes the current row 17 (DAIRIS SEGOVIA PEREZ / periodo 2405)
#    down to row 42. We then restore the "middle row" formatting (same as
#    row 16) on the freshly inserted rows by pasting formats only, so we
#    don't leave half-formatted blank rows behind.
# ---------------------------------------------------------------------------
$ws.Range("B17:J41").Insert()
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J41").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Fill in the 25 new worker rows (periodo mora 2507)
# ---------------------------------------------------------------------------
$newRows = @(
  @(17,  "33108217",   "ELEIDIS MARIA MARTINEZ PANZA",       1423500),
  @(18,  "22790985",   "LENIS DEL ROSARIO AGUILAR RAMOS",     1423500),
  @(19,  "1002321069", "ISAURA ACOSTA VILLAR",                1423500),
  @(20,  "45689297",   "NURIS PEREZ HERAZO",                  1000000),
  @(21,  "22793524",   "MARIA DEL CARMEN ARROYO ZUÑIGA",      1423500),
  @(22,  "45492177",   "LILIANA ROMERO PAZ",                  1423500),
  @(23,  "1100011464", "ATALIA MILAGRO BOLIVAR NAVARRO",      1423500),
  @(24,  "45756898",   "JUANITA DEL CARMEN ESPITIA MORALES",  1423500),
  @(25,  "73142180",   "HENRY RAMON HERRERA PLAZA",           1423500),
  @(26,  "1002390353", "KEYLEN SANCHEZ CUESTA",               1423500),
  @(27,  "26176370",   "MARIA VICTORIA CASTELLANOS GONZALEZ", 1423500),
  @(28,  "73150132",   "OSVALDO JOSE ROYERO CORONADO",        1423500),
  @(29,  "1081921984", "JESUS ALBERTO BERRIO RESTREPO",       1423500),
  @(30,  "45757694",   "BERTILDA GARCIA FORTICH",             1423500),
  @(31,  "1143385360", "MARIA DEL PILAR SEHUANES VERGARA",    1423500),
  @(32,  "45502601",   "ENITH MERCEDES NIETO BLANCO",         1423500),
  @(33,  "50982864",   "ANGELA CRISTINA LOPEZ GOMEZ",         1423500),
  @(34,  "64476123",   "MARLENE AYLEN CORREA SIERRA",         1000000),
  @(35,  "1143361235", "MARIA ALEJANDRA HERRERA HERRERA",     1423500),
  @(36,  "1047365006", "JUAN GUILLERMO MUÑOZ HERRERA",        1423500),
  @(37,  "1047452474", "EDGAR JOSE REALES ALVARES",           1423500),
  @(38,  "1043645205", "BENICIA PAOLA SEHUANES VERGARA",      1423500),
  @(39,  "1047504980", "LILIANA AGAMEZ ZUÑIGA",               1423500),
  @(40,  "1043296034", "ZHARICK PAOLA CAMPILLO CUADRO",       1423500),
  @(41,  "1047485489", "DAIRIS SEGOVIA PEREZ",                1423500)
)

foreach ($row in $newRows) {
  $r = $row[0]
  $ws.Cells.Item($r, 2).Value = "CC"
  $ws.Cells.Item($r, 3).Value = $row[1]
  $ws.Cells.Item($r, 4).Value = $row[2]
  $ws.Cells.Item($r, 5).Value = "2507"
  $ws.Cells.Item($r, 6).Value = 56940
  $ws.Cells.Item($r, 7).Value = $row[3]
}

# ---------------------------------------------------------------------------
# 3) Row 42 (former row 17) keeps periodo 2405 for DAIRIS SEGOVIA PEREZ, but
#    the "Valor Mora" (column G) is updated.
# ---------------------------------------------------------------------------
$ws.Range("G42").Value = 1423500

# ---------------------------------------------------------------------------
# 4) Update the summary header values
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 1497633
$ws.Range("C13").Value = 26
$ws.Range("F13").Value = 3
